# Natmi following Dr Hou advice
#
# Regenerate the LR-pair table for Ntn1-Unc5b: a new "ECs" sending/target
# cluster joins the existing "FAPs" and "sCs" clusters, expanding the 2x2
# cluster-pair grid (4 data rows) into a 3x3 grid (9 data rows), with freshly
# recomputed NATMI expression/specificity statistics for every cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 9,20

# Row 2: ECs -> ECs (Ntn1/Unc5b)
$data[0,0] = 'ECs'
$data[0,1] = 'Ntn1'
$data[0,2] = 'Unc5b'
$data[0,3] = 'ECs'
$data[0,4] = 2
$data[0,5] = 0.6666666666666666
$data[0,6] = 0.9305633333333333
$data[0,7] = 2.79169
$data[0,8] = 0.01768777137856805
$data[0,9] = 0.01768777137856806
$data[0,10] = 2
$data[0,11] = 0.6666666666666666
$data[0,12] = 5.108824666666666
$data[0,13] = 15.326474
$data[0,14] = 0.5049726372337502
$data[0,15] = 0.5049726372337501
$data[0,16] = 4.754084911228889
$data[0,17] = 42.78676420106
$data[0,18] = 0.008931840559823155
$data[0,19] = 0.008931840559823155

# Row 3: ECs -> FAPs (Ntn1/Unc5b)
$data[1,0] = 'ECs'
$data[1,1] = 'Ntn1'
$data[1,2] = 'Unc5b'
$data[1,3] = 'FAPs'
$data[1,4] = 2
$data[1,5] = 0.6666666666666666
$data[1,6] = 0.9305633333333333
$data[1,7] = 2.79169
$data[1,8] = 0.01768777137856805
$data[1,9] = 0.01768777137856806
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 3.403844333333333
$data[1,13] = 10.211533
$data[1,14] = 0.3364469054793339
$data[1,15] = 0.3364469054793339
$data[1,16] = 3.167492728974444
$data[1,17] = 28.50743456077
$data[1,18] = 0.005950995945145154
$data[1,19] = 0.005950995945145154

# Row 4: ECs -> sCs (Ntn1/Unc5b)
$data[2,0] = 'ECs'
$data[2,1] = 'Ntn1'
$data[2,2] = 'Unc5b'
$data[2,3] = 'sCs'
$data[2,4] = 2
$data[2,5] = 0.6666666666666666
$data[2,6] = 0.9305633333333333
$data[2,7] = 2.79169
$data[2,8] = 0.01768777137856805
$data[2,9] = 0.01768777137856806
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 1.604363666666667
$data[2,13] = 4.813091
$data[2,14] = 0.158580457286916
$data[2,15] = 0.1585804572869159
$data[2,16] = 1.492962001532222
$data[2,17] = 13.43665801379
$data[2,18] = 0.002804934873599746
$data[2,19] = 0.002804934873599746

# Row 5: FAPs -> ECs (Ntn1/Unc5b)
$data[3,0] = 'FAPs'
$data[3,1] = 'Ntn1'
$data[3,2] = 'Unc5b'
$data[3,3] = 'ECs'
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 44.154177
$data[3,7] = 132.462531
$data[3,8] = 0.8392647337471152
$data[3,9] = 0.8392647337471153
$data[3,10] = 2
$data[3,11] = 0.6666666666666666
$data[3,12] = 5.108824666666666
$data[3,13] = 15.326474
$data[3,14] = 0.5049726372337502
$data[3,15] = 0.5049726372337501
$data[3,16] = 225.575948593966
$data[3,17] = 2030.183537345694
$data[3,18] = 0.4238057259375619
$data[3,19] = 0.4238057259375619

# Row 6: FAPs -> FAPs (Ntn1/Unc5b)
$data[4,0] = 'FAPs'
$data[4,1] = 'Ntn1'
$data[4,2] = 'Unc5b'
$data[4,3] = 'FAPs'
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 44.154177
$data[4,7] = 132.462531
$data[4,8] = 0.8392647337471152
$data[4,9] = 0.8392647337471153
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 3.403844333333333
$data[4,13] = 10.211533
$data[4,14] = 0.3364469054793339
$data[4,15] = 0.3364469054793339
$data[4,16] = 150.293945174447
$data[4,17] = 1352.645506570023
$data[4,18] = 0.282368022547154
$data[4,19] = 0.282368022547154

# Row 7: FAPs -> sCs (Ntn1/Unc5b)
$data[5,0] = 'FAPs'
$data[5,1] = 'Ntn1'
$data[5,2] = 'Unc5b'
$data[5,3] = 'sCs'
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 44.154177
$data[5,7] = 132.462531
$data[5,8] = 0.8392647337471152
$data[5,9] = 0.8392647337471153
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 1.604363666666667
$data[5,13] = 4.813091
$data[5,14] = 0.158580457286916
$data[5,15] = 0.1585804572869159
$data[5,16] = 70.83935731036901
$data[5,17] = 637.5542157933211
$data[5,18] = 0.1330909852623993
$data[5,19] = 0.1330909852623993

# Row 8: sCs -> ECs (Ntn1/Unc5b)
$data[6,0] = 'sCs'
$data[6,1] = 'Ntn1'
$data[6,2] = 'Unc5b'
$data[6,3] = 'ECs'
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 7.525807
$data[6,7] = 22.577421
$data[6,8] = 0.1430474948743168
$data[6,9] = 0.1430474948743168
$data[6,10] = 2
$data[6,11] = 0.6666666666666666
$data[6,12] = 5.108824666666666
$data[6,13] = 15.326474
$data[6,14] = 0.5049726372337502
$data[6,15] = 0.5049726372337501
$data[6,16] = 38.44802843817266
$data[6,17] = 346.032255943554
$data[6,18] = 0.07223507073636511
$data[6,19] = 0.07223507073636509

# Row 9: sCs -> FAPs (Ntn1/Unc5b)
$data[7,0] = 'sCs'
$data[7,1] = 'Ntn1'
$data[7,2] = 'Unc5b'
$data[7,3] = 'FAPs'
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 7.525807
$data[7,7] = 22.577421
$data[7,8] = 0.1430474948743168
$data[7,9] = 0.1430474948743168
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 3.403844333333333
$data[7,13] = 10.211533
$data[7,14] = 0.3364469054793339
$data[7,15] = 0.3364469054793339
$data[7,16] = 25.61667551071033
$data[7,17] = 230.550079596393
$data[7,18] = 0.04812788698703476
$data[7,19] = 0.04812788698703475

# Row 10: sCs -> sCs (Ntn1/Unc5b)
$data[8,0] = 'sCs'
$data[8,1] = 'Ntn1'
$data[8,2] = 'Unc5b'
$data[8,3] = 'sCs'
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 7.525807
$data[8,7] = 22.577421
$data[8,8] = 0.1430474948743168
$data[8,9] = 0.1430474948743168
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 1.604363666666667
$data[8,13] = 4.813091
$data[8,14] = 0.158580457286916
$data[8,15] = 0.1585804572869159
$data[8,16] = 12.07413131314567
$data[8,17] = 108.667181818311
$data[8,18] = 0.02268453715091692
$data[8,19] = 0.02268453715091692

# Write the full 20-column block in one shot, replacing the old A2:T5 block
# and extending the used range to A1:T10.
$ws.Range("A2:T10").Value = $data
